$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 73, shifting the existing rows 73-115 down to 74-116
$ws.Rows.Item(73).Insert()

# Populate the newly inserted row 73 with the new weekly data point
$ws.Cells.Item(73, 1).Value = 6
$ws.Cells.Item(73, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(73, 3).Value = "Metropolitana"
$ws.Cells.Item(73, 4).Value = 44488
$ws.Cells.Item(73, 5).Value = 13
$ws.Cells.Item(73, 6).Value = 100112029
$ws.Cells.Item(73, 7).Value = "Orégano"
$ws.Cells.Item(73, 8).Value = "Sin especificar"
$ws.Cells.Item(73, 9).Value = "Primera"
$ws.Cells.Item(73, 10).Value = 34
$ws.Cells.Item(73, 11).Value = 8500
$ws.Cells.Item(73, 12).Value = 9000
$ws.Cells.Item(73, 13).Value = 8735
$ws.Cells.Item(73, 14).Value = "$/docena de atados"
$ws.Cells.Item(73, 15).Value = "Región Metropolitana"
$ws.Cells.Item(73, 16).Value = 2912
$ws.Cells.Item(73, 17).Value = 3
$ws.Cells.Item(73, 18).Value = "Hortaliza"
